# Auto-generated edit script applying cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '68.550.90'
$ws.Range('E2').Value = '  +1.36%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.955.01'
$ws.Range('E3').Value = '  +2.47%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '482.04'
$ws.Range('E5').Value = '  +5.55%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '150.18'
$ws.Range('E6').Value = '  +2.32%  '
$ws.Range('E7').Value = '  -0.27%  '
$ws.Range('E8').Value = '  -0.16%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.729'
$ws.Range('E9').Value = '  -2.37%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.169'
$ws.Range('E10').Value = '  +8.64%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0000355'
$ws.Range('E11').Value = '  +11.25%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '42.89'
$ws.Range('E12').Value = '  -1.99%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.593.00'
$ws.Range('E13').Value = '  +2.61%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '10.51'
$ws.Range('E14').Value = '  +1.23%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.002.52'
$ws.Range('E15').Value = '  +3.93%  '
$ws.Range('B16').Value = 'Uniswap'
$ws.Range('C16').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '14.80'
$ws.Range('E16').Value = '  -0.57%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '19.88'
$ws.Range('E18').Value = '  -1.17%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '1.13'
$ws.Range('E19').Value = '  -2.81%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '68.711.09'
$ws.Range('E20').Value = '  +1.48%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '436.88'
$ws.Range('E21').Value = '  +2.06%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '3.41'
$ws.Range('E22').Value = '  +4.72%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '14.51'
$ws.Range('E23').Value = '  -2.45%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '87.66'
$ws.Range('E24').Value = '  +0.95%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '11.03'
$ws.Range('E25').Value = '  +11.71%  '
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.57'
$ws.Range('E26').Value = '  +1.70%  '
$ws.Range('B27').Value = 'RenderToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.57'
$ws.Range('E27').Value = '  +2.04%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '38.46'
$ws.Range('E28').Value = '  +2.39%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.87'
$ws.Range('E29').Value = '  +6.83%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '716.55'
$ws.Range('E30').Value = '  -3.40%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '13.34'
$ws.Range('E31').Value = '  -3.50%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.129'
$ws.Range('E32').Value = '  -4.07%  '
$ws.Range('E33').Value = '  +3.60%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0₃0906'
$ws.Range('E34').Value = '  +31.92%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '41.99'
$ws.Range('E35').Value = '  -2.94%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '58.92'
$ws.Range('E36').Value = '  +3.00%  '
$ws.Range('E37').Value = '  -6.20%  '
$ws.Range('E38').Value = '  +0.04%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.998'
$ws.Range('E39').Value = '  -0.16%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.87'
$ws.Range('E40').Value = '  +6.72%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0471'
$ws.Range('E41').Value = '  -1.42%  '
$ws.Range('E42').Value = '  +10.22%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.01'
$ws.Range('E43').Value = '  -0.61%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.345'
$ws.Range('E44').Value = '  -3.13%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.141'
$ws.Range('E45').Value = '  +0.52%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.00'
$ws.Range('E46').Value = '  +0.01%  '
$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.18'
$ws.Range('E47').Value = '  +2.11%  '
$ws.Range('B48').Value = 'LidoDAOToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '3.44'
$ws.Range('E48').Value = '  -0.68%  '
$ws.Range('E49').Value = '  -1.06%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '146.48'
$ws.Range('E50').Value = '  +1.07%  '
$ws.Range('E51').Value = '  -1.31%  '
